# Auto-generated edit script
# Applies per-cell numeric updates to match the target diff across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 1520
$ws.Range("I54").Value = 1520
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1520
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1034
$ws.Range("N54").ClearContents()
$ws.Range("H132").Value = 2544.5732
$ws.Range("I132").Value = 1203.7167
$ws.Range("J132").Value = 7908
$ws.Range("K132").Value = 3611.1501
$ws.Range("L132").Value = 23724
$ws.Range("M132").Value = -1081.1501
$ws.Range("N132").Value = -28784
$ws.Range("H137").Value = 23123.45
$ws.Range("I137").Value = 27608.8
$ws.Range("J137").Value = 3188.5557
$ws.Range("K137").Value = 82826.39999999999
$ws.Range("L137").Value = 9565.667099999999
$ws.Range("M137").Value = -80276.39999999999
$ws.Range("N137").Value = -14665.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 846.92
$ws.Range("I97").Value = 955.5
$ws.Range("J97").Value = 412.6
$ws.Range("K97").Value = 955.5
$ws.Range("L97").Value = 412.6
$ws.Range("M97").Value = -459.5
$ws.Range("N97").Value = -1404.6
$ws.Range("H102").Value = 3257.1
$ws.Range("I102").Value = 1970
$ws.Range("K102").Value = 1970
$ws.Range("M102").Value = -348

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1670.9615
$ws.Range("I20").Value = 1654.3846
$ws.Range("J20").Value = 1687.5385
$ws.Range("K20").Value = 1654.3846
$ws.Range("L20").Value = 1687.5385
$ws.Range("M20").Value = -1407.3846
$ws.Range("N20").Value = -2181.5385
$ws.Range("H86").Value = 10438.333
$ws.Range("I86").Value = 7522.8887
$ws.Range("J86").Value = 19184.666
$ws.Range("K86").Value = 7522.8887
$ws.Range("L86").Value = 19184.666
$ws.Range("M86").Value = -6399.8887
$ws.Range("N86").Value = -21430.666
$ws.Range("H89").Value = 10438.333
$ws.Range("I89").Value = 7522.8887
$ws.Range("J89").Value = 19184.666
$ws.Range("K89").Value = 37614.4435
$ws.Range("L89").Value = 95923.33
$ws.Range("M89").Value = -31998.4435
$ws.Range("N89").Value = -107155.33
$ws.Range("H94").Value = 1480.8182
$ws.Range("I94").Value = 944.44446
$ws.Range("K94").Value = 944.44446
$ws.Range("M94").Value = -493.44446
$ws.Range("H99").Value = 3514.0952
$ws.Range("I99").Value = 5836.091
$ws.Range("J99").Value = 959.9
$ws.Range("K99").Value = 5836.091
$ws.Range("L99").Value = 959.9
$ws.Range("M99").Value = -4338.091
$ws.Range("N99").Value = -3955.9
$ws.Range("H105").Value = 2423.2703
$ws.Range("I105").Value = 2138.6365
$ws.Range("J105").Value = 2840.7334
$ws.Range("K105").Value = 2138.6365
$ws.Range("L105").Value = 2840.7334
$ws.Range("M105").Value = -391.6365000000001
$ws.Range("N105").Value = -6334.7334
$ws.Range("H134").Value = 22367.96
$ws.Range("I134").Value = 29498.756
$ws.Range("J134").Value = 4778.6665
$ws.Range("K134").Value = 88496.26800000001
$ws.Range("L134").Value = 14335.9995
$ws.Range("M134").Value = -85961.26800000001
$ws.Range("N134").Value = -19405.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2927106
$ws.Range("I62").Value = 6175402.5
$ws.Range("J62").Value = 3639.4
$ws.Range("K62").Value = 6175402.5
$ws.Range("L62").Value = 3639.4
$ws.Range("M62").Value = -6174778.5
$ws.Range("N62").Value = -4887.4
$ws.Range("H65").Value = 2927106
$ws.Range("I65").Value = 6175402.5
$ws.Range("J65").Value = 3639.4
$ws.Range("K65").Value = 30877012.5
$ws.Range("L65").Value = 18197
$ws.Range("M65").Value = -30873892.5
$ws.Range("N65").Value = -24437
$ws.Range("H105").Value = 635.03125
$ws.Range("I105").Value = 611.53845
$ws.Range("J105").Value = 736.8333
$ws.Range("K105").Value = 611.53845
$ws.Range("L105").Value = 736.8333
$ws.Range("M105").Value = 1135.46155
$ws.Range("N105").Value = -4230.8333
$ws.Range("H132").Value = 2031.1428
$ws.Range("I132").Value = 1166.25
$ws.Range("K132").Value = 3498.75
$ws.Range("M132").Value = -968.75
$ws.Range("H134").Value = 1694.6216
$ws.Range("I134").Value = 1062.238
$ws.Range("J134").Value = 2524.625
$ws.Range("K134").Value = 3186.714
$ws.Range("L134").Value = 7573.875
$ws.Range("M134").Value = -651.7139999999999
$ws.Range("N134").Value = -12643.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 14188
$ws.Range("I140").Value = 1167.5555
$ws.Range("J140").Value = 30928.572
$ws.Range("K140").Value = 3502.6665
$ws.Range("L140").Value = 92785.716
$ws.Range("M140").Value = 1677.3335
$ws.Range("N140").Value = -103145.716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5168.5186
$ws.Range("I80").Value = 5952.5
$ws.Range("J80").Value = 2928.5715
$ws.Range("K80").Value = 5952.5
$ws.Range("L80").Value = 2928.5715
$ws.Range("M80").Value = -4954.5
$ws.Range("N80").Value = -4924.5715
$ws.Range("H83").Value = 5168.5186
$ws.Range("I83").Value = 5952.5
$ws.Range("J83").Value = 2928.5715
$ws.Range("K83").Value = 29762.5
$ws.Range("L83").Value = 14642.8575
$ws.Range("M83").Value = -24770.5
$ws.Range("N83").Value = -24626.8575
$ws.Range("H97").Value = 1225.3636
$ws.Range("I97").Value = 971.5833
$ws.Range("K97").Value = 971.5833
$ws.Range("M97").Value = -475.5833
$ws.Range("H126").Value = 3707.111
$ws.Range("I126").Value = 3783.1667
$ws.Range("J126").Value = 3555
$ws.Range("K126").Value = 11349.5001
$ws.Range("L126").Value = 10665
$ws.Range("M126").Value = -8879.500100000001
$ws.Range("N126").Value = -15605
$ws.Range("H132").Value = 4667.0625
$ws.Range("I132").Value = 4921.2
$ws.Range("K132").Value = 14763.6
$ws.Range("M132").Value = -12233.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3452.842
$ws.Range("I68").Value = 1935
$ws.Range("J68").Value = 4556.727
$ws.Range("K68").Value = 1935
$ws.Range("L68").Value = 4556.727
$ws.Range("M68").Value = -1186
$ws.Range("N68").Value = -6054.727
$ws.Range("H71").Value = 3452.842
$ws.Range("I71").Value = 1935
$ws.Range("J71").Value = 4556.727
$ws.Range("K71").Value = 9675
$ws.Range("L71").Value = 22783.635
$ws.Range("M71").Value = -5931
$ws.Range("N71").Value = -30271.635
$ws.Range("H93").Value = 1458.1786
$ws.Range("I93").Value = 1415.2778
$ws.Range("J93").Value = 1535.4
$ws.Range("K93").Value = 1415.2778
$ws.Range("L93").Value = 1535.4
$ws.Range("M93").Value = -167.2778000000001
$ws.Range("N93").Value = -4031.4
$ws.Range("H100").Value = 43480810
$ws.Range("I100").Value = 2438.111
$ws.Range("J100").Value = 71431190
$ws.Range("K100").Value = 2438.111
$ws.Range("L100").Value = 71431190
$ws.Range("M100").Value = -1897.111
$ws.Range("N100").Value = -71432272
$ws.Range("H122").Value = 2694.7727
$ws.Range("I122").Value = 2361
$ws.Range("J122").Value = 2885.5
$ws.Range("K122").Value = 7083
$ws.Range("L122").Value = 8656.5
$ws.Range("M122").Value = -4633
$ws.Range("N122").Value = -13556.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4343.381
$ws.Range("J62").Value = 4437.4375
$ws.Range("L62").Value = 4437.4375
$ws.Range("N62").Value = -5685.4375
$ws.Range("H65").Value = 4343.381
$ws.Range("J65").Value = 4437.4375
$ws.Range("L65").Value = 22187.1875
$ws.Range("N65").Value = -28427.1875
$ws.Range("H81").Value = 2092.3809
$ws.Range("I81").Value = 878.3333
$ws.Range("J81").Value = 3711.111
$ws.Range("K81").Value = 1756.6666
$ws.Range("L81").Value = 7422.222
$ws.Range("M81").Value = -695.6666
$ws.Range("N81").Value = -9544.222
$ws.Range("H84").Value = 2092.3809
$ws.Range("I84").Value = 878.3333
$ws.Range("J84").Value = 3711.111
$ws.Range("K84").Value = 8783.333000000001
$ws.Range("L84").Value = 37111.11
$ws.Range("M84").Value = -3479.333000000001
$ws.Range("N84").Value = -47719.11
$ws.Range("H96").Value = 6488.6
$ws.Range("I96").Value = 2500
$ws.Range("K96").Value = 2500
$ws.Range("M96").Value = -1127
$ws.Range("H132").Value = 23018.604
$ws.Range("I132").Value = 46607.5
$ws.Range("J132").Value = 3058.7693
$ws.Range("K132").Value = 139822.5
$ws.Range("L132").Value = 9176.3079
$ws.Range("M132").Value = -137292.5
$ws.Range("N132").Value = -14236.3079
